# Updated cryptos list on Mon May  8 13:07:15 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row with
# newly scraped figures. Rows 42/43 additionally swap which coin
# (FraxShare / Algorand) occupies which rank, so their Coin (B) and Link (C)
# columns are rewritten too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2;  D="27.847.15";    E="  -3.84%  "},
    @{Row=3;  D="1.861.73";     E="  -2.74%  "},
    @{Row=4;  D="1.004";        E="  +0.13%  "},
    @{Row=5;  D="317.39";       E="  -2.29%  "},
    @{Row=6;  D="1.003";        E="  +0.25%  "},
    @{Row=7;  D="0.4367"},
    @{Row=8;  E="  -3.09%  "},
    @{Row=9;  D="0.07494";      E="  -2.61%  "},
    @{Row=10; D="0.9379";       E="  -4.28%  "},
    @{Row=11; D="21.29";        E="  -4.11%  "},
    @{Row=12; D="1.852.98";     E="  -3.00%  "},
    @{Row=13; D="6.743";        E="  -3.06%  "},
    @{Row=14; D="5.454";        E="  -4.09%  "},
    @{Row=15; D="0.06844";      E="  -2.52%  "},
    @{Row=16; D="1.004";        E="  +0.08%  "},
    @{Row=17; D="81.58";        E="  -2.89%  "},
    @{Row=18; D="0.000009072";  E="  -3.81%  "},
    @{Row=19; D="1.003";        E="  +0.25%  "},
    @{Row=20; D="15.92";        E="  -4.25%  "},
    @{Row=21; D="27.854.88";    E="  -3.82%  "},
    @{Row=22; D="5.117";        E="  -3.78%  "},
    @{Row=23; D="11.07";        E="  +1.26%  "},
    @{Row=24; D="2.107.80";     E="  -1.12%  "},
    @{Row=25; D="2.007";        E="  -4.10%  "},
    @{Row=26; D="154.07"},
    @{Row=27; D="18.38";        E="  -3.31%  "},
    @{Row=28; D="5.469";        E="  -3.88%  "},
    @{Row=29; D="113.09";       E="  -3.75%  "},
    @{Row=30; D="1.717";        E="  -8.03%  "},
    @{Row=31; D="0.09022";      E="  -2.89%  "},
    @{Row=32; D="0.8127";       E="  -6.04%  "},
    @{Row=33; D="4.817";        E="  -5.56%  "},
    @{Row=34; D="1.174";        E="  -5.76%  "},
    @{Row=35; E="  -3.40%  "},
    @{Row=36; E="  +0.27%  "},
    @{Row=37; D="1.119";        E="  -3.23%  "},
    @{Row=38; D="0.05483";      E="  -3.89%  "},
    @{Row=39; D="0.01975";      E="  -3.15%  "},
    @{Row=40; D="2.953";        E="  -1.34%  "},
    @{Row=41; D="0.5261";       E="  -4.40%  "},
    @{Row=42; B="Algorand";  C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.1706";  E="  -2.53%  "},
    @{Row=43; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs";      D="7.028";   E="  -6.14%  "},
    @{Row=44; D="8.783";        E="  -6.30%  "},
    @{Row=45; D="0.06763";      E="  -1.93%  "},
    @{Row=46; D="0.4904";       E="  -5.08%  "},
    @{Row=47; D="10.63";        E="  -5.18%  "},
    @{Row=48; D="107.50";       E="  -2.49%  "},
    @{Row=49; D="1.680";        E="  -5.61%  "},
    @{Row=50; D="1.002";        E="  +0.15%  "},
    @{Row=51; D="1.884";        E="  -14.12%  "}
)

foreach ($item in $updates) {
    $r = $item.Row

    if ($item.ContainsKey("B")) {
        $ws.Cells.Item($r, 2).Value = $item.B
    }
    if ($item.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $item.C
    }
    if ($item.ContainsKey("D")) {
        # Force text: the Price column stores figures like "1.004" or
        # "27.847.15" as plain text, not numbers, so mark the cell as Text
        # before assigning or Excel will silently coerce numeric-looking
        # strings into real numbers.
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $item.E
    }
}
